$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from
# 45204 (2023-10-05) to 45205 (2023-10-06) for every data row (2-180).
$ws.Range("C2:C180").Value = 45205
